# "break out stock.yaml completed"
# - Fix bsecode values in E60:E63 which were stored as text; store as numbers.
# - Append the newly completed breakout rows 64-66 (06:45:32 scan), keeping
#   their bsecode column (E) as text, matching the source data feed format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Convert E60:E63 bsecode values from text to numeric
$ws.Range("E60").Value = 509930
$ws.Range("E61").Value = 590024
$ws.Range("E62").Value = 543220
$ws.Range("E63").Value = 532155

# Row 64
$ws.Range("A64").Value = "26/06/2024 06:45:32"
$ws.Range("B64").Value = 1
$ws.Range("C64").Value = "SUPREMEIND"
$ws.Range("D64").Value = "Supreme Industries Limited"
$ws.Range("E64").Value = "'509930"
$ws.Range("F64").Value = -1.59
$ws.Range("G64").Value = 5794.6
$ws.Range("H64").Value = 42571

# Row 65
$ws.Range("A65").Value = "26/06/2024 06:45:32"
$ws.Range("B65").Value = 2
$ws.Range("C65").Value = "FACT"
$ws.Range("D65").Value = "Fertilizers And Chemicals Travancore Limited"
$ws.Range("E65").Value = "'590024"
$ws.Range("F65").Value = 2.89
$ws.Range("G65").Value = 1027.95
$ws.Range("H65").Value = 1990905

# Row 66
$ws.Range("A66").Value = "26/06/2024 06:45:32"
$ws.Range("B66").Value = 3
$ws.Range("C66").Value = "MAXHEALTH"
$ws.Range("D66").Value = "Max Healthcare Institute Ltd"
$ws.Range("E66").Value = "'543220"
$ws.Range("F66").Value = -2.07
$ws.Range("G66").Value = 875
$ws.Range("H66").Value = 1099451
